# Reorders the data rows (2-12) of the historical-distance table so that
# each row's (title, timestamp, uri) triple lands in a new row position,
# while the "historical distance" / "time bucket" columns stay "unknown".
# This mirrors the source edit: the underlying dataset/time-bucket analysis
# JSON was regenerated, which reshuffled row order without altering any
# individual record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state per row: Title (col A), Timestamp (col B), Uri (col E).
# Columns C (historical distance) and D (time bucket) remain "unknown".
$titles = @{
    2  = "Presidential Ratings"
    3  = "New Hampshire: Election Tools, Deadlines, Dates, Rules, and Links"
    4  = "Sanders Still Leads in NH as Primary Nears, Poll Finds"
    5  = "Find Your Local League"
    6  = "2020 President - Sabato's Crystal Ball"
    7  = "Biden Leads In New Hampshire Poll"
    8  = "EXCLUSIVE: N.H. Voters Oppose Medicare For All by 2-to-1 Margin"
    9  = "Center for Public Opinion"
    10 = "Post NHJournal Poll Results"
    11 = "2020 Election Forecast"
    12 = "2020 Electoral Interactive Map"
}

$timestamps = @{
    2  = "1-01-01T00:00:00UTC"
    3  = "1-01-01T00:00:00UTC"
    4  = "2020-02-02T17:15:03UTC"
    5  = "1-01-01T00:00:00UTC"
    6  = "1-01-01T00:00:00UTC"
    7  = "2020-09-29T21:30:00UTC"
    8  = "2020-02-07T13:54:55UTC"
    9  = "1-01-01T00:00:00UTC"
    10 = "2001-08-07T00:00:00UTC"
    11 = "2020-08-12T06:30:00UTC"
    12 = "1-01-01T00:00:00UTC"
}

$uris = @{
    2  = "https://insideelections.com/ratings/president"
    3  = "https://www.vote.org/state/new-hampshire/"
    4  = "https://www.nbcboston.com/news/politics/nh-primary-poll-sanders-still-leads/2070807/"
    5  = "https://www.lwv.org/local-leagues/find-local-league"
    6  = "http://centerforpolitics.org/crystalball/2020-president/"
    7  = "https://amgreatness.com/2020/09/29/biden-leads-in-new-hampshire-poll/"
    8  = "https://www.insidesources.com/n-h-voters-oppose-medicare-for-all-by-2-to-1-margin/"
    9  = "https://www.uml.edu/Research/public-opinion/polls/2020/NH-Sept.aspx"
    10 = "https://docs.google.com/document/d/1Z1lbFrjnIBmYZxWLZkzo_1nVg3uTRLcerSSauihlLdM/edit"
    11 = "https://projects.fivethirtyeight.com/2020-election-forecast/"
    12 = "https://abcnews.go.com/Politics/2020-Electoral-Interactive-Map?basemap=71662160&promoref=brandpromo"
}

for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = $titles[$r]
    $ws.Cells.Item($r, 2).Value = $timestamps[$r]
    $ws.Cells.Item($r, 5).Value = $uris[$r]
}
